# Reorder the "Recorded By" names in column G for the
# "Session Analysis Results" sheet.
#
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, system, System"   -> "system, backup@backdoor.com, System"
#   "admin@admin.com, dnasr281@gmail.com"   -> "dnasr281@gmail.com, admin@admin.com"
#
# Every cell in column G whose text matches one of the patterns above is
# rewritten with the reordered value; everything else is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value = "system, backup@backdoor.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
